# Apply updated crypto market data (price + 1h volume change) per latest GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 49 and 50 swap which coin (Coin name + Link) they describe, and all four rows
# 49/50 get refreshed Price / Volume(1h) figures as part of this data update.

# --- Plain text / percentage / non-numeric-looking price updates ---
$plainUpdates = @{
    'D2' = '69.681.41'
    'E2' = '  +0.57%  '
    'D3' = '3.688.32'
    'E3' = '  +0.43%  '
    'E4' = '  -0.02%  '
    'E5' = '  -1.17%  '
    'E6' = '  +1.12%  '
    'E7' = '  -0.04%  '
    'E8' = '  +1.29%  '
    'E9' = '  -0.11%  '
    'E10' = '  +3.26%  '
    'E11' = '  +1.52%  '
    'E12' = '  +0.97%  '
    'E13' = '  +1.60%  '
    'D14' = '3.693.05'
    'E14' = '  +0.48%  '
    'D15' = '69.656.01'
    'E15' = '  +0.61%  '
    'E16' = '  +2.47%  '
    'E17' = '  +0.78%  '
    'E18' = '  +0.64%  '
    'E19' = '  +0.67%  '
    'E20' = '  -2.04%  '
    'E21' = '  -0.20%  '
    'E22' = '  +0.07%  '
    'D23' = '3.833.88'
    'E23' = '  +0.40%  '
    'E24' = '  +4.77%  '
    'E26' = '  +0.71%  '
    'E27' = '  +0.05%  '
    'E28' = '  +0.03%  '
    'E29' = '  -2.56%  '
    'E30' = '  +1.04%  '
    'E31' = '  +0.12%  '
    'E32' = '  +2.08%  '
    'E33' = '  -0.53%  '
    'E34' = '  -2.13%  '
    'D35' = '3.674.04'
    'E35' = '  +0.22%  '
    'E36' = '  +3.34%  '
    'E37' = '  -1.88%  '
    'E39' = '  +1.10%  '
    'E40' = '  -0.06%  '
    'E41' = '  +2.08%  '
    'E42' = '  +0.68%  '
    'E43' = '  -0.71%  '
    'E44' = '  -1.31%  '
    'E45' = '  +2.12%  '
    'E46' = '  -2.11%  '
    'E47' = '  -0.12%  '
    'E48' = '  -1.90%  '
    'B49' = 'SuiNetwork'
    'C49' = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
    'E49' = '  +0.24%  '
    'B50' = 'Cosmos'
    'C50' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E50' = '  +1.02%  '
    'E51' = '  -0.19%  '
}
foreach ($ref in $plainUpdates.Keys) {
    $ws.Range($ref).Value = $plainUpdates[$ref]
}

# --- Price updates that look like plain numbers (e.g. '1.00', '0.0000128') ---
# These must be forced to Text format, otherwise Excel auto-converts them to
# numeric values (losing trailing zeros / switching to scientific notation),
# since column D is stored as text in this workbook.
$textForcedUpdates = @{
    'D4' = '1.00'
    'D5' = '667.70'
    'D6' = '160.21'
    'D13' = '32.88'
    'D17' = '16.16'
    'D19' = '470.37'
    'D20' = '9.78'
    'D22' = '79.81'
    'D24' = '0.0000128'
    'D26' = '10.96'
    'D27' = '9.07'
    'D30' = '2.01'
    'D32' = '0.166'
    'D33' = '26.74'
    'D34' = '6.48'
    'D36' = '8.47'
    'D37' = '6.13'
    'D39' = '2.26'
    'D41' = '176.99'
    'D44' = '47.01'
    'D46' = '27.68'
    'D49' = '1.07'
    'D50' = '7.86'
}
foreach ($ref in $textForcedUpdates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = '@'
    $cell.Value = $textForcedUpdates[$ref]
}

